# Normalize the "Recorded By" (column G) name ordering on the
# "Session Analysis Results" sheet so that the synced value matches the
# canonical ordering from the main repo.
#
# The sync only reorders the comma-separated recorder names within a cell;
# it never adds/removes/changes a recorder. Observed reorderings:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"                 -> "admin@admin.com, System"
#   "admin@admin.com, dnasr281@gmail.com"     -> "dnasr281@gmail.com, admin@admin.com"
#   "backup@backdoor.com, System, system"     -> "backup@backdoor.com, system, System"
# Any other combination (e.g. "backup@backdoor.com, System") is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "backup@backdoor.com, System, system" = "backup@backdoor.com, system, System"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($null -eq $current) { continue }
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
